$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "... have one or more Contact methods ..."
#         -> "... have zero or more Contact methods ..."
# ------------------------------------------------------------------

# Locate the whole target sentence first, so we only touch this occurrence.
$sentence1 = $d.Content
$sentence1.Find.ClearFormatting()
$found1 = $sentence1.Find.Execute( `
    "A User must be able to have one or more Contact methods, which are also existence-dependent on the User.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    # Narrow down to just the word "one" inside that sentence.
    $word1 = $d.Range($sentence1.Start, $sentence1.End)
    $word1.Find.ClearFormatting()
    $word1.Find.Execute("one", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Replace the word's text.
    $word1.Text = "zero"

    # Nudge a character-formatting property on just this word (and revert it)
    # so the run containing it gets split out from its neighbours, matching
    # how Word itself splits a run when only part of it is edited/selected.
    $word1.Bold = 1
    $word1.Bold = 0
}

# ------------------------------------------------------------------
# Edit 2: "The system must allow a User to exist without a Role assigned."
#         -> "The system must not allow a User to exist without a Role assigned."
# ------------------------------------------------------------------

$sentence2 = $d.Content
$sentence2.Find.ClearFormatting()
$found2 = $sentence2.Find.Execute( `
    "The system must allow a User to exist without a Role assigned.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $word2 = $d.Range($sentence2.Start, $sentence2.End)
    $word2.Find.ClearFormatting()
    $word2.Find.Execute("allow", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

    # Insert the new word right before "allow"; $word2 then covers "not allow".
    $word2.InsertBefore("not ")

    # Narrow to just the inserted "not " text and nudge formatting on it so
    # it becomes its own run, split from the following "allow ..." text.
    $notRange = $d.Range($word2.Start, $word2.Start + 4)
    $notRange.Bold = 1
    $notRange.Bold = 0
}
